$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into the "Price" (D) column as literal text (coinranking
# price strings look numeric, e.g. "522.61", and Excel's Value setter would
# otherwise silently coerce them to a Number cell and drop the text formatting).
# A leading apostrophe forces text entry; re-applying the "Normal" style
# afterwards clears the quote-prefix formatting flag Excel adds so the cell's
# style index is left untouched (matches the original file, which carries no
# explicit style on these cells).
function Set-PriceText($row, $value) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

function Set-VolumeText($row, $value) {
    $ws.Cells.Item($row, 5).Value = $value
}

# row -> new Price (D) text (rows not listed here keep their original Price)
$prices = @{
    2  = "57.231.29"
    3  = "3.070.99"
    5  = "522.61"
    6  = "135.56"
    7  = "0.999"
    8  = "3.071.33"
    9  = "0.471"
    14 = "3.599.69"
    15 = "25.08"
    17 = "57.257.99"
    18 = "3.062.18"
    20 = "12.39"
    21 = "7.79"
    22 = "347.07"
    23 = "1.00"
    24 = "68.84"
    25 = "0.497"
    31 = "1.84"
    32 = "5.82"
    33 = "20.99"
    34 = "158.33"
    35 = "4.78"
    38 = "25.42"
    40 = "0.0655"
    42 = "4.00"
    43 = "0.690"
    44 = "2.410.72"
    45 = "36.64"
    47 = "3.109.41"
    48 = "0.0258"
    49 = "5.96"
    50 = "0.934"
}

# row -> new Volume(1h) (E) text (all rows get an E update)
$volumes = @{
    2  = "  -2.08%  "
    3  = "  -2.07%  "
    4  = "  +0.02%  "
    5  = "  -2.07%  "
    6  = "  -5.34%  "
    7  = "  -0.05%  "
    8  = "  -2.03%  "
    9  = "  +4.45%  "
    10 = "  +0.84%  "
    11 = "  -3.44%  "
    12 = "  +1.16%  "
    13 = "  +1.38%  "
    14 = "  -2.10%  "
    15 = "  -2.30%  "
    16 = "  -4.43%  "
    17 = "  -2.09%  "
    18 = "  -2.50%  "
    19 = "  -4.53%  "
    20 = "  -3.86%  "
    21 = "  -2.49%  "
    22 = "  +1.49%  "
    23 = "  +0.01%  "
    24 = "  +1.47%  "
    25 = "  -3.22%  "
    28 = "  -9.96%  "
    29 = "  +0.01%  "
    30 = "  -5.52%  "
    31 = "  -3.33%  "
    32 = "  -10.10%  "
    33 = "  -0.74%  "
    34 = "  +0.23%  "
    35 = "  -0.16%  "
    36 = "  -7.42%  "
    38 = "  -2.89%  "
    39 = "  -4.12%  "
    40 = "  -2.61%  "
    41 = "  -6.40%  "
    42 = "  -0.47%  "
    43 = "  -2.59%  "
    44 = "  +4.61%  "
    45 = "  +0.13%  "
    46 = "  +0.07%  "
    47 = "  -2.05%  "
    48 = "  -2.55%  "
    49 = "  -2.28%  "
    50 = "  -7.29%  "
    51 = "  -6.70%  "
}

foreach ($row in $prices.Keys) {
    Set-PriceText $row $prices[$row]
}
foreach ($row in $volumes.Keys) {
    Set-VolumeText $row $volumes[$row]
}

# PEPE (row 28) price uses a Unicode subscript-3 digit ("0.0\u20830840" -> "...0838").
# Concatenating that codepoint into a plain numeric-looking literal makes Excel's
# parser misread the whole string as a number, so build it explicitly and use
# the same text-forcing technique as the other Price cells.
$subscript3 = [string][char]0x2083
$pepePrice = "0.0" + $subscript3 + "0838"
Set-PriceText 28 $pepePrice

# Rows 26 and 27 swap coin identities: Kaspa moves up to row 26 and
# Binance-PegBSC-USD moves down to row 27 (their Price/Volume figures also change).
$ws.Cells.Item(26, 2).Value = "Kaspa"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-PriceText 26 "0.165"
Set-VolumeText 26 "  -2.88%  "

$ws.Cells.Item(27, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-PriceText 27 "0.999"
Set-VolumeText 27 "  -0.13%  "
